$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.341.30"
$ws.Range("D3").Value = "3.687.31"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "680.77"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "159.09"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "7.13"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("D13").Value = "4.309.49"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "32.42"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "3.679.30"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "69.304.40"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "468.56"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "80.01"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "3.835.43"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -4.93%  "
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  -3.72%  "
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "26.91"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").Value = "3.675.61"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("E36").Value = "  -4.94%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "6.25"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("D43").Value = "170.66"
$ws.Range("E43").Value = "  +4.28%  "
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").Value = "28.23"
$ws.Range("E46").Value = "  -5.18%  "
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "2.68"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").Value = "0.000275"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("E51").Value = "  -2.85%  "
